$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (bold, border, centered) from existing header cell H1
# so the new header cells I1/J1 reuse the same style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row values
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I (I0) and J (IF) columns, rows 2-13
$values = @(
    @(8, 9),
    @(7, 8),
    @(8, 9),
    @(1, 4),
    @(8, 9),
    @(8, 9),
    @(7, 8),
    @(9, 9),
    @(1, 4),
    @(1, 3),
    @(2, 3),
    @(1, 1)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
